$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D that render as plain decimal numbers need a leading
# apostrophe so Excel keeps storing them as text (matching the source data,
# which uses text-formatted price strings like "26.916.79" or "1.001").

$ws.Range("D2").Value = "26.916.79"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.809.46"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'310.41"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4618"
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("D8").Value = "'0.3711"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").Value = "'0.07366"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'0.8756"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "1.854.90"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "'5.361"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'92.22"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "'6.509"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").Value = "'0.07043"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'0.000008705"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "26.902.35"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'5.325"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").Value = "'10.64"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").Value = "2.023.70"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "'1.894"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "'151.31"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "'2.154"
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'115.87"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "'0.08907"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'0.7570"
$ws.Range("E32").Value = "  -4.60%  "
$ws.Range("D33").Value = "'1.157"
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("D34").Value = "'4.455"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'2.916"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "'0.01970"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'0.05249"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'2.423"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "'0.5321"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'7.215"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'8.517"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").Value = "'0.4985"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "'10.37"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9999"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.671"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'103.74"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").Value = "'0.06291"
$ws.Range("E51").Value = "  -1.57%  "
